$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.141.48"
$ws.Range("E2").Value = "  -1.71%  "
$ws.Range("D3").Value = "'3.745.03"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'614.81"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").Value = "'176.98"
$ws.Range("E6").Value = "  -0.72%  "
$ws.Range("D7").Value = "'3.744.82"
$ws.Range("E7").Value = "  -1.48%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.527"
$ws.Range("E9").Value = "  -1.81%  "
$ws.Range("E10").Value = "  -2.47%  "
$ws.Range("D11").Value = "'6.55"
$ws.Range("E11").Value = "  +3.26%  "
$ws.Range("D12").Value = "'0.484"
$ws.Range("E12").Value = "  -2.23%  "
$ws.Range("E14").Value = "  -1.74%  "
$ws.Range("D15").Value = "'4.367.48"
$ws.Range("E15").Value = "  -1.60%  "
$ws.Range("D16").Value = "'3.750.17"
$ws.Range("E16").Value = "  -1.39%  "
$ws.Range("D17").Value = "'69.226.46"
$ws.Range("E17").Value = "  -1.51%  "
$ws.Range("E18").Value = "  -2.85%  "
$ws.Range("D19").Value = "'7.46"
$ws.Range("E19").Value = "  -1.97%  "
$ws.Range("D20").Value = "'16.32"
$ws.Range("E20").Value = "  -2.18%  "
$ws.Range("D21").Value = "'498.62"
$ws.Range("E21").Value = "  -3.37%  "
$ws.Range("D22").Value = "'9.35"
$ws.Range("E22").Value = "  -2.91%  "
$ws.Range("D23").Value = "'0.722"
$ws.Range("E23").Value = "  -0.92%  "
$ws.Range("D24").Value = "'2.53"
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D25").Value = "'85.82"
$ws.Range("E25").Value = "  -2.73%  "
$ws.Range("D26").Value = "'12.82"
$ws.Range("E26").Value = "  -3.59%  "
$ws.Range("D27").Value = "'10.73"
$ws.Range("E27").Value = "  -3.66%  "
$ws.Range("E28").Value = "  -3.35%  "
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("D30").Value = "'2.51"
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("D31").Value = "'2.94"
$ws.Range("E31").Value = "  +3.49%  "
$ws.Range("E32").Value = "  +1.75%  "
$ws.Range("D33").Value = "'30.51"
$ws.Range("E33").Value = "  -4.29%  "
$ws.Range("E34").Value = "  -1.37%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("D37").Value = "'6.09"
$ws.Range("E37").Value = "  -2.04%  "
$ws.Range("D38").Value = "'0.347"
$ws.Range("E38").Value = "  +1.97%  "
$ws.Range("E39").Value = "  +2.95%  "
$ws.Range("D40").Value = "'449.80"
$ws.Range("E40").Value = "  +6.30%  "
$ws.Range("E41").Value = "  -5.17%  "
$ws.Range("E42").Value = "  +9.13%  "
$ws.Range("D43").Value = "'49.72"
$ws.Range("E43").Value = "  -3.03%  "
$ws.Range("D44").Value = "'44.82"
$ws.Range("E44").Value = "  +1.58%  "
$ws.Range("E45").Value = "  -2.44%  "
$ws.Range("D46").Value = "'2.940.23"
$ws.Range("E46").Value = "  -4.20%  "
$ws.Range("D47").Value = "'0.0359"
$ws.Range("E47").Value = "  -1.86%  "
$ws.Range("D48").Value = "'27.29"
$ws.Range("E48").Value = "  -1.13%  "
$ws.Range("D50").Value = "'138.48"
$ws.Range("E50").Value = "  +1.84%  "
$ws.Range("E51").Value = "  -0.31%  "
